$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Venus Explorer (column F) values
$ws.Range("F3").Value = 800
$ws.Range("F4").Value = 100
$ws.Range("F6").Value = 0.9
$ws.Range("F10").Value = 70
$ws.Range("F16").Value = 10
$ws.Range("F19").Value = 0.25
$ws.Range("F21").Value = "BPSK_Viterbi"

# Move the active selection from E5 to F5
$ws.Activate()
$ws.Range("F5").Select()
